# Generate Report for Handoff
# Adds a new file (a635c615-cd83-4968-96a3-721b37eadc04.md) entry to all
# three report sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$commitUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5483604ad6d937853bf7e5f2c8bc854a1871263f/e2e/a635c615-cd83-4968-96a3-721b37eadc04.md"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Duplicate row 2 (with formatting) into the new row 3.
$wsOverview.Rows(2).Copy()
$wsOverview.Rows(3).Insert()

$wsOverview.Range("A3").Value = "a635c615-cd83-4968-96a3-721b37eadc04.md"
$wsOverview.Range("B3").Value = "e2e\a635c615-cd83-4968-96a3-721b37eadc04.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-24 14:44:40"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $commitUrl, [type]::Missing, [type]::Missing, "e2e\a635c615-cd83-4968-96a3-721b37eadc04.md") | Out-Null
$wsOverview.Range("B3").Font.Underline = $true
$wsOverview.Range("B3").Font.Color = 15570276

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Rows(2).Copy()
$wsZhCn.Rows(3).Insert()

$wsZhCn.Range("A3").Value = "a635c615-cd83-4968-96a3-721b37eadc04.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = "a635c615-cd83-4968-96a3-721b37eadc04.3ec1a4d60513514c9c1d82d0bdb125047a061ca8.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-24 14:44:35"
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("M3").Value = "False"
$wsZhCn.Range("O3").Value = "False"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $commitUrl, [type]::Missing, [type]::Missing, "a635c615-cd83-4968-96a3-721b37eadc04.md") | Out-Null
$wsZhCn.Range("A3").Font.Underline = $true
$wsZhCn.Range("A3").Font.Color = 15570276

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Rows(2).Copy()
$wsDeDe.Rows(3).Insert()

$wsDeDe.Range("A3").Value = "a635c615-cd83-4968-96a3-721b37eadc04.md"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = "a635c615-cd83-4968-96a3-721b37eadc04.3ec1a4d60513514c9c1d82d0bdb125047a061ca8.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-24 14:44:40"
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("M3").Value = "False"
$wsDeDe.Range("O3").Value = "False"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $commitUrl, [type]::Missing, [type]::Missing, "a635c615-cd83-4968-96a3-721b37eadc04.md") | Out-Null
$wsDeDe.Range("A3").Font.Underline = $true
$wsDeDe.Range("A3").Font.Color = 15570276

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P3"))
